# Weekly CompStat (13th Precinct) refresh: bump the report volume/number and
# the covered date range, and replace last week's crime-complaint figures
# with the newly collected numbers for rows 15-28 and 31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- helpers -----------------------------------------------------------
# Plain numeric value; optionally copy the number format/style from another
# cell so cells that are switching from a "N/A" text placeholder back to a
# real number pick up the right style (integer / percent / etc.).
function Set-NumCell($sheet, $row, $col, $value, $fmtRow = $null, $fmtCol = $null) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.Value = $value
    if ($fmtRow -ne $null) {
        $sheet.Cells.Item($fmtRow, $fmtCol).Copy()
        $cell.PasteSpecial($xlPasteFormats)
    }
}

# Shared "N/A" style text value ("0" or "***.*"). Forces the cell to be
# stored as text, then borrows the exact style from a cell that already
# carries that text so the style id matches (General-format text style).
function Set-TextCell($sheet, $row, $col, $text, $fmtRow, $fmtCol) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $sheet.Cells.Item($fmtRow, $fmtCol).Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

# --- header: volume/number + date range ---------------------------------
$ws.Range("A8").Value = "Volume 31   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/18/2024  Through  11/24/2024"

# --- row 15 (Rape) -------------------------------------------------------
Set-NumCell  $ws 15 6 5
Set-TextCell $ws 15 7 "0"     15 4
Set-TextCell $ws 15 8 "***.*" 15 5
Set-NumCell  $ws 15 9 22
Set-NumCell  $ws 15 11 100
Set-NumCell  $ws 15 12 29.411764705882
Set-NumCell  $ws 15 13 100
Set-NumCell  $ws 15 14 22.222222222222

# --- row 16 (Robbery) -----------------------------------------------------
Set-NumCell $ws 16 3 4
Set-NumCell $ws 16 4 2
Set-NumCell $ws 16 5 100
Set-NumCell $ws 16 6 29
Set-NumCell $ws 16 7 9
Set-NumCell $ws 16 8 222.222222222222
Set-NumCell $ws 16 9 159
Set-NumCell $ws 16 10 158
Set-NumCell $ws 16 11 0.632911392405
Set-NumCell $ws 16 12 -10.674157303370
Set-NumCell $ws 16 13 -4.216867469879
Set-NumCell $ws 16 14 -83.385579937304

# --- row 17 (Fel. Assault) -------------------------------------------------
Set-NumCell $ws 17 3 3
Set-NumCell $ws 17 4 2
Set-NumCell $ws 17 5 50
Set-NumCell $ws 17 6 17
Set-NumCell $ws 17 7 14
Set-NumCell $ws 17 8 21.428571428571
Set-NumCell $ws 17 9 212
Set-NumCell $ws 17 10 207
Set-NumCell $ws 17 11 2.415458937198
Set-NumCell $ws 17 12 -3.636363636363
Set-NumCell $ws 17 13 69.6
Set-NumCell $ws 17 14 -20.599250936329

# --- row 18 (Burglary) -----------------------------------------------------
Set-NumCell $ws 18 3 4
Set-NumCell $ws 18 4 2
Set-NumCell $ws 18 5 100
Set-NumCell $ws 18 6 17
Set-NumCell $ws 18 8 13.333333333333
Set-NumCell $ws 18 9 182
Set-NumCell $ws 18 10 200
Set-NumCell $ws 18 11 -9
Set-NumCell $ws 18 12 -28.90625
Set-NumCell $ws 18 13 -20.175438596491
Set-NumCell $ws 18 14 -87.619047619047

# --- row 19 (Gr. Larceny) ---------------------------------------------------
Set-NumCell $ws 19 3 24
Set-NumCell $ws 19 4 22
Set-NumCell $ws 19 5 9.090909090909
Set-NumCell $ws 19 6 120
Set-NumCell $ws 19 7 89
Set-NumCell $ws 19 8 34.831460674157
Set-NumCell $ws 19 9 976
Set-NumCell $ws 19 10 950
Set-NumCell $ws 19 11 2.736842105263
Set-NumCell $ws 19 12 3.499469777306
Set-NumCell $ws 19 13 -26.395173453997
Set-NumCell $ws 19 14 -60.96

# --- row 20 (G.L.A.) ---------------------------------------------------------
Set-TextCell $ws 20 3 "0" 23 3
Set-NumCell  $ws 20 4 2
Set-NumCell  $ws 20 5 -100
Set-NumCell  $ws 20 6 2
Set-NumCell  $ws 20 7 4
Set-NumCell  $ws 20 8 -50
Set-NumCell  $ws 20 9 38
Set-NumCell  $ws 20 10 58
Set-NumCell  $ws 20 11 -34.482758620689
Set-NumCell  $ws 20 12 -42.424242424242
Set-NumCell  $ws 20 13 -24
Set-NumCell  $ws 20 14 -96.729776247848

# --- row 21 (TOTAL) -----------------------------------------------------------
Set-NumCell $ws 21 3 36
Set-NumCell $ws 21 4 30
Set-NumCell $ws 21 5 20
Set-NumCell $ws 21 6 190
Set-NumCell $ws 21 8 45.038167938931
Set-NumCell $ws 21 9 1591
Set-NumCell $ws 21 10 1585
Set-NumCell $ws 21 11 0.378548895899
Set-NumCell $ws 21 12 -5.410225921522
Set-NumCell $ws 21 13 -16.657936092194
Set-NumCell $ws 21 14 -75.062695924764

# --- row 22 (Transit) -----------------------------------------------------
Set-TextCell $ws 22 3 "0" 23 3
Set-NumCell  $ws 22 5 -100
Set-NumCell  $ws 22 6 5
Set-NumCell  $ws 22 8 0
Set-NumCell  $ws 22 10 82
Set-NumCell  $ws 22 11 -29.268292682926
Set-NumCell  $ws 22 12 -30.120481927710
Set-NumCell  $ws 22 13 -14.705882352941

# --- row 23 (Housing) -----------------------------------------------------
Set-TextCell $ws 23 4 "0"     23 3
Set-TextCell $ws 23 5 "***.*" 22 14
Set-TextCell $ws 23 6 "0"     23 3
Set-NumCell  $ws 23 7 1
Set-NumCell  $ws 23 8 -100

# --- row 24 (Petit Larceny) -----------------------------------------------
Set-NumCell $ws 24 3 98
Set-NumCell $ws 24 4 66
Set-NumCell $ws 24 5 48.484848484848
Set-NumCell $ws 24 6 312
Set-NumCell $ws 24 7 233
Set-NumCell $ws 24 8 33.905579399141
Set-NumCell $ws 24 9 2757
Set-NumCell $ws 24 10 2043
Set-NumCell $ws 24 11 34.948604992657
Set-NumCell $ws 24 12 29.924599434495
Set-NumCell $ws 24 13 62.176470588235

# --- row 25 (Retail Theft) -------------------------------------------------
Set-NumCell $ws 25 3 94
Set-NumCell $ws 25 4 45
Set-NumCell $ws 25 5 108.888888888889
Set-NumCell $ws 25 6 278
Set-NumCell $ws 25 7 168
Set-NumCell $ws 25 8 65.476190476190
Set-NumCell $ws 25 9 2351
Set-NumCell $ws 25 10 1512
Set-NumCell $ws 25 11 55.489417989418
Set-NumCell $ws 25 12 46.571072319202

# --- row 26 (Misd. Assault) -------------------------------------------------
Set-NumCell $ws 26 3 13
Set-NumCell $ws 26 4 4
Set-NumCell $ws 26 5 225
Set-NumCell $ws 26 6 62
Set-NumCell $ws 26 7 38
Set-NumCell $ws 26 8 63.157894736842
Set-NumCell $ws 26 9 562
Set-NumCell $ws 26 10 458
Set-NumCell $ws 26 11 22.707423580786
Set-NumCell $ws 26 12 20.342612419700
Set-NumCell $ws 26 13 46.736292428198

# --- row 27 (UCR Rape*) ----------------------------------------------------
Set-NumCell $ws 27 3 2
Set-NumCell $ws 27 7 1
Set-NumCell $ws 27 8 500
Set-NumCell $ws 27 9 35
Set-NumCell $ws 27 11 59.090909090909
Set-NumCell $ws 27 12 25

# --- row 28 (Other Sex Crimes) ----------------------------------------------
Set-NumCell  $ws 28 3 1
Set-TextCell $ws 28 4 "0"     27 4
Set-TextCell $ws 28 5 "***.*" 27 5
Set-NumCell  $ws 28 7 7
Set-NumCell  $ws 28 8 0

# --- row 31 (Hate Crimes) ----------------------------------------------------
Set-NumCell $ws 31 3 1 31 6
Set-NumCell $ws 31 4 1 31 6
Set-NumCell $ws 31 5 0 31 8
Set-NumCell $ws 31 6 1
Set-NumCell $ws 31 7 5
Set-NumCell $ws 31 8 -80
Set-NumCell $ws 31 9 10
Set-NumCell $ws 31 10 12
Set-NumCell $ws 31 11 -16.666666666666
Set-NumCell $ws 31 12 -54.545454545454

$ws.Application.CutCopyMode = $false
